$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the current row 254 (shifts existing rows 254-325 down to 255-326)
$ws.Rows.Item(254).Insert()

# Populate the newly inserted row 254 with the new weekly record
$ws.Cells.Item(254, 1).Value  = 8
$ws.Cells.Item(254, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(254, 3).Value  = "Coquimbo"
$ws.Cells.Item(254, 4).Value  = 44932
$ws.Cells.Item(254, 5).Value  = 4
$ws.Cells.Item(254, 6).Value  = 100112031
$ws.Cells.Item(254, 7).Value  = "Poroto verde"
$ws.Cells.Item(254, 8).Value  = "Magnum"
$ws.Cells.Item(254, 9).Value  = "Primera"
$ws.Cells.Item(254, 10).Value = 400
$ws.Cells.Item(254, 11).Value = 23000
$ws.Cells.Item(254, 12).Value = 24000
$ws.Cells.Item(254, 13).Value = 23500
$ws.Cells.Item(254, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(254, 15).Value = "Provincia de Limar" + [char]0x00ED
$ws.Cells.Item(254, 16).Value = 940
$ws.Cells.Item(254, 17).Value = 25
$ws.Cells.Item(254, 18).Value = "Hortaliza"
